$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-10 Monday" "2025-03-11 Tuesday"
Replace-Text "815×4=" "647×7="
Replace-Text "174×9=" "443×6="
Replace-Text "413×8=" "116×2="
Replace-Text "625×9=" "137×7="
Replace-Text "706×9=" "799×4="
Replace-Text "571×5=" "776×6="
Replace-Text "521×2=" "191×7="
Replace-Text "739×2=" "621×2="
Replace-Text "741×9=" "352×8="
Replace-Text "880×3=" "731×6="
Replace-Text "285×3=" "336×9="
Replace-Text "835×8=" "177×9="
Replace-Text "211×8=" "257×9="
Replace-Text "327×5=" "535×6="
Replace-Text "597×7=" "674×4="
Replace-Text "599×7=" "976×5="
Replace-Text "159×5=" "303×2="
Replace-Text "531×6=" "213×8="
Replace-Text "413×7=" "871×6="
Replace-Text "743×5=" "544×4="
Replace-Text "963×4=" "855×4="
Replace-Text "673×7=" "303×9="
Replace-Text "904×2=" "857×9="
Replace-Text "133×8=" "175×5="
Replace-Text "735×7=" "429×7="
